$wb = $excel.ActiveWorkbook

# --- Rename sheets ---
$wb.Worksheets.Item(1).Name = "GNG_TO-16504778114935884"
$wb.Worksheets.Item(2).Name = "NB_TO-1650477813861685"
$wb.Worksheets.Item(3).Name = "RS_TO-16504778138626602"
$wb.Worksheets.Item(4).Name = "TOL_TO-1650477813924661"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16504778139856637"

# --- Sheet 1: GNG_TO ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16504778114525573.csv"
$ws1.Range("B3").Value = "GNG_stims-16504778114765916.csv"
$ws1.Range("B4").Value = "go_stims-16504778114775524.csv"
$ws1.Range("B5").Value = "GNG_stims-16504778114925544.csv"

# --- Sheet 2: NB_TO ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "OB-16504778130225558.csv"
$ws2.Range("B3").Value = "TB-1650477813841696.csv"
$ws2.Range("B4").Value = "TB-16504778137675872.csv"
$ws2.Range("B5").Value = "ZB-match_1-16504778125755532.csv"
$ws2.Range("B6").Value = "OB-1650477812779585.csv"
$ws2.Range("B7").Value = "ZB-match_9-1650477812652555.csv"
$ws2.Range("B8").Value = "ZB-match_3-16504778127135904.csv"
$ws2.Range("B9").Value = "OB-16504778131585526.csv"
$ws2.Range("B10").Value = "TB-16504778137115867.csv"

# --- Sheet 3: RS_TO --- (no cell value changes)

# --- Sheet 4: TOL_TO ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16504778138776953.csv"
$ws4.Range("B3").Value = "ZM_stims-16504778138646638.csv"
$ws4.Range("B4").Value = "MM_stims-16504778139086642.csv"
$ws4.Range("B5").Value = "ZM_stims-16504778138786597.csv"
$ws4.Range("B6").Value = "MM_stims-1650477813924661.csv"
$ws4.Range("B7").Value = "ZM_stims-16504778139096618.csv"

# --- Sheet 5: vSAT_TO ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-1650477813970661.csv"
$ws5.Range("B3").Value = "SAT_stims-16504778139286602.csv"
$ws5.Range("B4").Value = "SAT_stims-16504778139396648.csv"
$ws5.Range("B5").Value = "vSAT_stims-16504778139546614.csv"
